# Daily TGP (Terminal Gate Price) update: shift effective dates forward one
# day and refresh the Diesel/ULP/PULP/e10 prices for each terminal row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = 46071
$ws.Range("D8").Value = 157.63
$ws.Range("E8").Value = 148.41999999999999
$ws.Range("F8").Value = 158.41999999999999
$ws.Range("G8").Value = 148.31

# Row 9
$ws.Range("A9").Value = 46071
$ws.Range("D9").Value = 157.63
$ws.Range("E9").Value = 148.41999999999999
$ws.Range("F9").Value = 158.41999999999999
$ws.Range("G9").Value = 148.31

# Row 10
$ws.Range("A10").Value = 46071
$ws.Range("D10").Value = 159.27000000000001
$ws.Range("E10").Value = 151.13
$ws.Range("F10").Value = 161.13
$ws.Range("G10").Value = 151.37

# Row 11
$ws.Range("A11").Value = 46070
$ws.Range("D11").Value = 158.33000000000001
$ws.Range("E11").Value = 149.19
$ws.Range("F11").Value = 159.19
$ws.Range("G11").Value = 149.08000000000001

# Row 12
$ws.Range("A12").Value = 46070
$ws.Range("D12").Value = 158.33000000000001
$ws.Range("E12").Value = 149.19
$ws.Range("F12").Value = 159.19
$ws.Range("G12").Value = 149.08000000000001

# Row 13
$ws.Range("A13").Value = 46070
$ws.Range("D13").Value = 159.71
$ws.Range("E13").Value = 151.87
$ws.Range("F13").Value = 161.87
$ws.Range("G13").Value = 152.11000000000001

# Row 17
$ws.Range("A17").Value = 46071
$ws.Range("D17").Value = 163.65
$ws.Range("E17").Value = 154.69
$ws.Range("F17").Value = 164.69

# Row 18
$ws.Range("A18").Value = 46070
$ws.Range("D18").Value = 164.07
$ws.Range("E18").Value = 155.41
$ws.Range("F18").Value = 165.41

# Row 22
$ws.Range("A22").Value = 46071
$ws.Range("D22").Value = 159.03
$ws.Range("E22").Value = 150.76
$ws.Range("F22").Value = 160.36000000000001
$ws.Range("G22").Value = 152.52000000000001

# Row 23
$ws.Range("A23").Value = 46071
$ws.Range("D23").Value = 164.24
$ws.Range("E23").Value = 156.83000000000001
$ws.Range("F23").Value = 166.83

# Row 24
$ws.Range("A24").Value = 46071
$ws.Range("D24").Value = 164.43
$ws.Range("E24").Value = 157.36000000000001
$ws.Range("F24").Value = 167.36

# Row 25
$ws.Range("A25").Value = 46071
$ws.Range("D25").Value = 164.44
$ws.Range("E25").Value = 156.87
$ws.Range("F25").Value = 166.87
$ws.Range("G25").Value = 157.72999999999999

# Row 26
$ws.Range("A26").Value = 46071
$ws.Range("D26").Value = 164.07
$ws.Range("E26").Value = 158.46
$ws.Range("F26").Value = 168.46

# Row 27
$ws.Range("A27").Value = 46070
$ws.Range("D27").Value = 159.51
$ws.Range("E27").Value = 151.53
$ws.Range("F27").Value = 161.13
$ws.Range("G27").Value = 153.29

# Row 28
$ws.Range("A28").Value = 46070
$ws.Range("D28").Value = 164.68
$ws.Range("E28").Value = 157.57
$ws.Range("F28").Value = 167.57

# Row 29
$ws.Range("A29").Value = 46070
$ws.Range("D29").Value = 164.87
$ws.Range("E29").Value = 158.11000000000001
$ws.Range("F29").Value = 168.11

# Row 30
$ws.Range("A30").Value = 46070
$ws.Range("D30").Value = 164.88
$ws.Range("E30").Value = 157.63
$ws.Range("F30").Value = 167.63
$ws.Range("G30").Value = 158.47999999999999

# Row 31
$ws.Range("A31").Value = 46070
$ws.Range("D31").Value = 164.5
$ws.Range("E31").Value = 159.22
$ws.Range("F31").Value = 169.22

# Row 35
$ws.Range("A35").Value = 46071
$ws.Range("D35").Value = 157.77000000000001
$ws.Range("E35").Value = 148.82
$ws.Range("F35").Value = 157.82

# Row 36
$ws.Range("A36").Value = 46070
$ws.Range("D36").Value = 158.22
$ws.Range("E36").Value = 149.56
$ws.Range("F36").Value = 158.56

# Row 40
$ws.Range("A40").Value = 46071
$ws.Range("D40").Value = 163.96
$ws.Range("E40").Value = 156.03
$ws.Range("F40").Value = 166.03

# Row 41
$ws.Range("A41").Value = 46071
$ws.Range("D41").Value = 163.68
$ws.Range("E41").Value = 156.46
$ws.Range("F41").Value = 166.46

# Row 42
$ws.Range("A42").Value = 46070
$ws.Range("D42").Value = 164.42
$ws.Range("E42").Value = 156.85
$ws.Range("F42").Value = 166.85

# Row 43
$ws.Range("A43").Value = 46070
$ws.Range("D43").Value = 164.14
$ws.Range("E43").Value = 157.27000000000001
$ws.Range("F43").Value = 167.27

# Row 47
$ws.Range("A47").Value = 46071
$ws.Range("D47").Value = 159.01
$ws.Range("E47").Value = 150.65
$ws.Range("F47").Value = 160.65

# Row 48
$ws.Range("A48").Value = 46071
$ws.Range("D48").Value = 158.71
$ws.Range("E48").Value = 150.63
$ws.Range("F48").Value = 160.63

# Row 49
$ws.Range("A49").Value = 46070
$ws.Range("D49").Value = 159.94
$ws.Range("E49").Value = 151.08000000000001
$ws.Range("F49").Value = 161.08000000000001

# Row 50
$ws.Range("A50").Value = 46070
$ws.Range("D50").Value = 159.63
$ws.Range("E50").Value = 151.06
$ws.Range("F50").Value = 161.06

# Row 54
$ws.Range("A54").Value = 46071
$ws.Range("D54").Value = 173.14
$ws.Range("E54").Value = 163.69999999999999
$ws.Range("F54").Value = 173.7

# Row 55
$ws.Range("A55").Value = 46071
$ws.Range("D55").Value = 162.52000000000001
$ws.Range("E55").Value = 162.72
$ws.Range("F55").Value = 172.72

# Row 56
$ws.Range("A56").Value = 46071
$ws.Range("D56").Value = 162.28

# Row 57
$ws.Range("A57").Value = 46071
$ws.Range("D57").Value = 163.18
$ws.Range("E57").Value = 157.13999999999999

# Row 58
$ws.Range("A58").Value = 46071
$ws.Range("D58").Value = 158.94999999999999
$ws.Range("E58").Value = 153.04
$ws.Range("F58").Value = 163.04

# Row 59
$ws.Range("A59").Value = 46071
$ws.Range("D59").Value = 166.09
$ws.Range("E59").Value = 162.25

# Row 60
$ws.Range("A60").Value = 46070
$ws.Range("D60").Value = 173.59
$ws.Range("E60").Value = 164.53
$ws.Range("F60").Value = 174.53

# Row 61
$ws.Range("A61").Value = 46070
$ws.Range("D61").Value = 162.96
$ws.Range("E61").Value = 163.22
$ws.Range("F61").Value = 173.22

# Row 62
$ws.Range("A62").Value = 46070
$ws.Range("D62").Value = 162.72999999999999

# Row 63
$ws.Range("A63").Value = 46070
$ws.Range("D63").Value = 163.59
$ws.Range("E63").Value = 157.63999999999999

# Row 64
$ws.Range("A64").Value = 46070
$ws.Range("D64").Value = 159.36000000000001
$ws.Range("E64").Value = 153.54
$ws.Range("F64").Value = 163.54

# Row 65
$ws.Range("A65").Value = 46070
$ws.Range("D65").Value = 166.5
$ws.Range("E65").Value = 163.05000000000001
